$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text in C1 (was "PrecioUnitarioendolares")
$ws.Range("C1").Value = "Precio Unitario en dolares"

# Update the current selection to match the edited file (rows 20:21 selected)
$ws.Range("A20:XFD21").Select()
